# Generate Report for Handoff
# Refresh the "Latest Handoff Date(time)" column for the files that were
# re-handed-off (rows 4-10 of each sheet) with the new handoff timestamps.

$wb = $excel.ActiveWorkbook

# Overview sheet: column D = "Latest Handoff Date"
$overview = $wb.Worksheets.Item("Overview")
for ($r = 4; $r -le 10; $r++) {
    $overview.Cells.Item($r, 4).Value = "2016-03-25 11:17:03"
}

# zh-cn sheet: column E = "Latest Handoff Datetime"
$zhcn = $wb.Worksheets.Item("zh-cn")
for ($r = 4; $r -le 10; $r++) {
    $zhcn.Cells.Item($r, 5).Value = "2016-03-25 11:16:53"
}

# de-de sheet: column E = "Latest Handoff Datetime"
$dede = $wb.Worksheets.Item("de-de")
for ($r = 4; $r -le 10; $r++) {
    $dede.Cells.Item($r, 5).Value = "2016-03-25 11:17:03"
}
